$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value = -12.584
$ws.Range("B9").Value = 6.636
$ws.Range("C11").Value = -12.934
$ws.Range("B18").Value = 5.972
$ws.Range("B20").Value = 6.37
$ws.Range("D21").Value = -7.805
